$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column 9 width: match column C's width (copy width) ---
$ws.Columns.Item(9).ColumnWidth = $ws.Columns.Item(3).ColumnWidth

# --- Row 1: J1 style change (remove fill, same border/font/alignment as H1/I1) ---
$j1 = $ws.Range("J1")
$j1.Font.Bold = $true
$j1.HorizontalAlignment = -4108
$j1.VerticalAlignment = -4160
$j1.Borders.Item(7).LineStyle = 1
$j1.Borders.Item(10).LineStyle = 1

# --- Row 2: add C2 ---
$ws.Range("C2").Value = "LXBC "

# --- Row 3: A3 value change 1 -> 3 ---
$ws.Range("A3").Value = 3

# --- Row 4: remove A4 ---
$ws.Range("A4").ClearContents()

# --- Row 7: new content (Bổn Mạng Nhà 2026 entry, year 2026) ---
$ws.Range("A7").Value = 2
$ws.Range("B7").Value = 2026
$ws.Range("C7").Value = "LXBC"
$ws.Range("D7").Value = "Ảnh"
$ws.Range("E7").Value = "Bổn Mạng Nhà 2026"
$ws.Range("I7").Value = "Đang Cập Nhập"

# --- Row 8: new row (Tư Liệu 2025 entry) ---
$ws.Range("A8").Value = 1
$ws.Range("C8").Value = "LXBC"
$ws.Range("D8").Value = "Ảnh"
$ws.Range("E8").Value = "Tư Liệu 2025"
$ws.Range("F8").Value = "gửi ảnh vào link để làm tư liệu bổn mạng "
$ws.Range("G8").Value = "https://drive.google.com/drive/folders/1IVawCMt9xO_6Cnvzh2S28Q6U66pePz7e"
$ws.Range("I8").Value = "Gửi Ảnh"
$ws.Range("H8").Value = "nhau.jpg"

# --- Selection state ---
$ws.Range("H9").Select()
